$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 0.7506253333333333
$ws.Range("I2").Value = 0.5569881999559233
$ws.Range("J2").Value = 0.5569881999559233
$ws.Range("M2").Value = 1.114648666666667
$ws.Range("N2").Value = 3.343946
$ws.Range("O2").Value = 0.4003231847851749
$ws.Range("P2").Value = 0.4003231847851748
$ws.Range("Q2").Value = 0.8366835269662223
$ws.Range("R2").Value = 7.530151742696
$ws.Range("S2").Value = 0.222975290094117
$ws.Range("T2").Value = 0.222975290094117

# Row 3
$ws.Range("G3").Value = 0.7506253333333333
$ws.Range("I3").Value = 0.5569881999559233
$ws.Range("J3").Value = 0.5569881999559233
$ws.Range("O3").Value = 0.3070172855255452
$ws.Range("P3").Value = 0.3070172855255452
$ws.Range("Q3").Value = 0.6416723164084444
$ws.Range("R3").Value = 5.775050847675999
$ws.Range("S3").Value = 0.1710050052202272
$ws.Range("T3").Value = 0.1710050052202272

# Row 4
$ws.Range("G4").Value = 0.7506253333333333
$ws.Range("I4").Value = 0.5569881999559233
$ws.Range("J4").Value = 0.5569881999559233
$ws.Range("M4").Value = 0.716238
$ws.Range("N4").Value = 2.148714
$ws.Range("O4").Value = 0.2572350246303296
$ws.Range("P4").Value = 0.2572350246303295
$ws.Range("Q4").Value = 0.537626387496
$ws.Range("R4").Value = 4.838637487463999
$ws.Range("S4").Value = 0.1432768733344649
$ws.Range("T4").Value = 0.1432768733344648

# Row 5
$ws.Range("G5").Value = 0.7506253333333333
$ws.Range("I5").Value = 0.5569881999559233
$ws.Range("J5").Value = 0.5569881999559233
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.09863499999999999
$ws.Range("N5").Value = 0.295905
$ws.Range("O5").Value = 0.03542450505895045
$ws.Range("P5").Value = 0.03542450505895045
$ws.Range("Q5").Value = 0.07403792975333331
$ws.Range("R5").Value = 0.6663413677799999
$ws.Range("S5").Value = 0.01973103130711431
$ws.Range("T5").Value = 0.01973103130711431

# Row 6
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 0.6666666666666666
$ws.Range("G6").Value = 0.597025
$ws.Range("H6").Value = 1.791075
$ws.Range("I6").Value = 0.4430118000440768
$ws.Range("J6").Value = 0.4430118000440768
$ws.Range("M6").Value = 1.114648666666667
$ws.Range("N6").Value = 3.343946
$ws.Range("O6").Value = 0.4003231847851749
$ws.Range("P6").Value = 0.4003231847851748
$ws.Range("Q6").Value = 0.6654731202166668
$ws.Range("R6").Value = 5.98925808195
$ws.Range("S6").Value = 0.1773478946910579
$ws.Range("T6").Value = 0.1773478946910579

# Row 7
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 0.6666666666666666
$ws.Range("G7").Value = 0.597025
$ws.Range("H7").Value = 1.791075
$ws.Range("I7").Value = 0.4430118000440768
$ws.Range("J7").Value = 0.4430118000440768
$ws.Range("O7").Value = 0.3070172855255452
$ws.Range("P7").Value = 0.3070172855255452
$ws.Range("Q7").Value = 0.5103670202583334
$ws.Range("R7").Value = 4.593303182324999
$ws.Range("S7").Value = 0.1360122803053181
$ws.Range("T7").Value = 0.136012280305318

# Row 8
$ws.Range("E8").Value = 2
$ws.Range("F8").Value = 0.6666666666666666
$ws.Range("G8").Value = 0.597025
$ws.Range("H8").Value = 1.791075
$ws.Range("I8").Value = 0.4430118000440768
$ws.Range("J8").Value = 0.4430118000440768
$ws.Range("M8").Value = 0.716238
$ws.Range("N8").Value = 2.148714
$ws.Range("O8").Value = 0.2572350246303296
$ws.Range("P8").Value = 0.2572350246303295
$ws.Range("Q8").Value = 0.42761199195
$ws.Range("R8").Value = 3.84850792755
$ws.Range("S8").Value = 0.1139581512958648
$ws.Range("T8").Value = 0.1139581512958647

# Row 9
$ws.Range("E9").Value = 2
$ws.Range("F9").Value = 0.6666666666666666
$ws.Range("G9").Value = 0.597025
$ws.Range("H9").Value = 1.791075
$ws.Range("I9").Value = 0.4430118000440768
$ws.Range("J9").Value = 0.4430118000440768
$ws.Range("K9").Value = 1
$ws.Range("L9").Value = 0.3333333333333333
$ws.Range("M9").Value = 0.09863499999999999
$ws.Range("N9").Value = 0.295905
$ws.Range("O9").Value = 0.03542450505895045
$ws.Range("P9").Value = 0.03542450505895045
$ws.Range("Q9").Value = 0.05888756087499999
$ws.Range("R9").Value = 0.529988047875
$ws.Range("S9").Value = 0.01569347375183614
$ws.Range("T9").Value = 0.01569347375183614

